$wb = $excel.ActiveWorkbook

# --- Fix segment and powertrain split values on the "Data" sheet ---
$data = $wb.Worksheets.Item("Data")

$data.Range("M2").Value = 0.7

$data.Range("F3").Value = 0.1
$data.Range("I3").Value = 0.8
$data.Range("J3").Value = 0.9
$data.Range("M3").Value = 0.3

$data.Range("F4").Value = 0.8
$data.Range("I4").Value = 0.1
$data.Range("J4").Value = 0

# --- Page setup (both sheets) ---
foreach ($ws in $wb.Worksheets) {
    $ws.PageSetup.PaperSize = 9
    $ws.PageSetup.Orientation = 1
}

# --- View state: "Data" becomes the active/selected sheet with G21 selected ---
$data.Activate()
$null = $data.Range("G21").Select()
